$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 52; existing rows 52-71 shift down to 53-72.
$ws.Rows("52:52").Insert()

# Preserve the date-formatted style used by column D in this table.
$ws.Range("D52").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Fill in the new row 52 with the new Puerro price record.
$ws.Range("A52").Value2 = 9
$ws.Range("B52").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C52").Value2 = "Metropolitana"
$ws.Range("D52").Value2 = 44489
$ws.Range("E52").Value2 = 13
$ws.Range("F52").Value2 = 100112005
$ws.Range("G52").Value2 = "Puerro"
$ws.Range("H52").Value2 = "Sin especificar"
$ws.Range("I52").Value2 = "Primera"
$ws.Range("J52").Value2 = 160
$ws.Range("K52").Value2 = 7000
$ws.Range("L52").Value2 = 8000
$ws.Range("M52").Value2 = 7500
$ws.Range("N52").Value2 = "$/paquete 20 unidades"
$ws.Range("O52").Value2 = "Provincia de Chacabuco"
$ws.Range("P52").Value2 = 375
$ws.Range("Q52").Value2 = 20
$ws.Range("R52").Value2 = "Hortaliza"
